$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# The localization report is re-generated: the "zh-cn" and "de-de" sheets
# each gain a populated "Latest Target File" (I6), "Latest Handback File"
# (J6), "Latest Handback DateTime" (K6) and "Error Detail" (P6) for the
# 2e340cf0-... row (row 6), plus a hyperlink on the new I6 cell. The
# "Error Detail" column (P) is also widened to fit the new text.
# ---------------------------------------------------------------------------

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen column P ("Error Detail") so the new message is readable.
    $ws.Columns.Item(16).ColumnWidth = 40

    # I6: "Latest Target File" now links to the handed-back markdown file.
    $i6 = $ws.Range("I6")
    $i6.Value = "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md"
    $i6.Style = "HyperLink"
    $ws.Hyperlinks.Add(
        $i6,
        "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa15af2ed472e40ac797dc924e5aa25893ddf519/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md",
        [Type]::Missing,
        [Type]::Missing,
        "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md"
    ) | Out-Null

    # J6: "Latest Handback File" - the xlf that was handed back.
    $ws.Range("J6").Value = "2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.151d020a437088accfe1576b6e5c5ed33cf41f2b.$sheetName.xlf"

    if ($sheetName -eq "zh-cn") {
        # zh-cn: the handback isn't the latest version -> surface the error
        # in the "Latest Handback DateTime" cell and record when the
        # (stale) handback happened in "Error Detail".
        $ws.Range("K6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aa15af2ed472e40ac797dc924e5aa25893ddf519/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e83ec9973823d5aed1f6cb65ff65711fa0b44cb2/e2e/2e340cf0-cfd9-47ef-a0a2-6096b1fa233c.md."
    }
    else {
        # de-de: handback matches the latest version -> plain timestamp.
        $ws.Range("K6").Value = "2016-08-30 16:51:53"
    }

    # P6: Error Detail / handback timestamp column.
    $ws.Range("P6").Value = "2016-08-30 16:51:29"
}
